# Edit script: restructure "Input" sheet columns/headers to the new standard
# template layout, add derived email columns, and drop the now-unused
# empty helper cell in "갑지" / "을지" (I4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$firstDataRow = 2
$lastDataRow  = 6
$stagingBase  = 29   # staging columns start at 30 (col 1 -> staging col 30)

# ---------------------------------------------------------------------------
# PHASE 1: stage a copy of every old column (1..17 == A..Q) for the data rows
# so later overwrites of columns 1..16 can't clobber a value we still need.
# ---------------------------------------------------------------------------
for ($oldCol = 1; $oldCol -le 17; $oldCol++) {
    $srcRange = $ws.Range($ws.Cells.Item($firstDataRow, $oldCol), $ws.Cells.Item($lastDataRow, $oldCol))
    $dstCell  = $ws.Cells.Item($firstDataRow, $stagingBase + $oldCol)
    $srcRange.Copy($dstCell)
}

function Stage($oldCol) {
    return $stagingBase + $oldCol
}

# ---------------------------------------------------------------------------
# PHASE 2: rebuild the data rows (2..6) in the new 16-column layout by
# pulling from the staged copies (preserves original cell typing, e.g. the
# date-like strings stay text instead of being re-parsed as serial dates).
# ---------------------------------------------------------------------------
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    # A: 발주일자 <- old C (발주일)
    $ws.Cells.Item($r, (Stage 3)).Copy($ws.Cells.Item($r, 1))
    # B: 납기일자 <- old D (납기일)
    $ws.Cells.Item($r, (Stage 4)).Copy($ws.Cells.Item($r, 2))
    # C: 거래처명 <- old A (거래처명)
    $ws.Cells.Item($r, (Stage 1)).Copy($ws.Cells.Item($r, 3))
    # E: 납품처명 <- old B (현장명)
    $ws.Cells.Item($r, (Stage 2)).Copy($ws.Cells.Item($r, 5))
    # G: 프로젝트명 <- old B (현장명)
    $ws.Cells.Item($r, (Stage 2)).Copy($ws.Cells.Item($r, 7))
    # H: 대분류 <- old N (대분류)
    $ws.Cells.Item($r, (Stage 14)).Copy($ws.Cells.Item($r, 8))
    # I: 중분류 <- old O (중분류)
    $ws.Cells.Item($r, (Stage 15)).Copy($ws.Cells.Item($r, 9))
    # J: 소분류 <- old P (소분류)
    $ws.Cells.Item($r, (Stage 16)).Copy($ws.Cells.Item($r, 10))
    # K: 품목명 <- old F (품목)
    $ws.Cells.Item($r, (Stage 6)).Copy($ws.Cells.Item($r, 11))
    # L: 규격 <- old G (규격)
    $ws.Cells.Item($r, (Stage 7)).Copy($ws.Cells.Item($r, 12))
    # M: 수량 <- old H (수량) [numeric]
    $ws.Cells.Item($r, (Stage 8)).Copy($ws.Cells.Item($r, 13))
    # N: 단가 <- old J (단가) [numeric]
    $ws.Cells.Item($r, (Stage 10)).Copy($ws.Cells.Item($r, 14))
    # O: 총금액 <- old M (합계) [numeric]
    $ws.Cells.Item($r, (Stage 13)).Copy($ws.Cells.Item($r, 15))
    # P: 비고 <- old Q (비고) (row 4 had none, leave blank)
    $ws.Cells.Item($r, (Stage 17)).Copy($ws.Cells.Item($r, 16))

    # D: 거래처 이메일 (new, derived from buyer name)
    $buyer = $ws.Cells.Item($r, 3).Text
    $ws.Cells.Item($r, 4).Value = "$buyer@example.com"
    # F: 납품처 이메일 (new, constant placeholder)
    $ws.Cells.Item($r, 6).Value = "delivery@example.com"
}

# row 4 has no "비고" in the new layout either - make sure it's empty
$ws.Cells.Item(4, 16).Value = ""

# ---------------------------------------------------------------------------
# PHASE 3: wipe the staging area now that the rebuild is done.
# ---------------------------------------------------------------------------
$stagingRange = $ws.Range($ws.Cells.Item($firstDataRow, $stagingBase + 1), $ws.Cells.Item($lastDataRow, $stagingBase + 17))
$stagingRange.Value = ""

# ---------------------------------------------------------------------------
# PHASE 4: drop the now-unused trailing column so the sheet dimension goes
# from A1:Q6 back down to A1:P6.
# ---------------------------------------------------------------------------
$ws.Columns.Item(17).Delete()

# ---------------------------------------------------------------------------
# PHASE 5: rewrite the header row with the new field names, and drop the
# bold/bordered header style (new template headers are plain).
# ---------------------------------------------------------------------------
$ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, 16)).Style = "Normal"

$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 갑지 / 을지: drop the stray empty helper cell at I4.
# ---------------------------------------------------------------------------
foreach ($sheetName in @("갑지", "을지")) {
    $s = $wb.Worksheets.Item($sheetName)
    $s.Cells.Item(4, 9).Value = ""
}
